# Update the "UnitMass" (column C) values in the two report tables on
# Sheet1 to reflect the refreshed catalog data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    5  = 71
    6  = 102
    7  = 24
    8  = 1
    9  = 56
    10 = 12
    11 = 13
    12 = 25
    13 = 74
    15 = 149
    16 = 14
    17 = 7
    18 = 175
    19 = 73
    20 = 38
    21 = 50
    27 = 68
    28 = 43
    30 = 39
    31 = 82
    34 = 15
    36 = 81
    37 = 61
    38 = 123
    39 = 19
    40 = 18
    41 = 126
    42 = 95
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}
